$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the header row (new rows 2 and 3);
# everything that was row 2.. shifts down to row 4..
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Copy the formatting (styles/borders) of the row immediately below (the old
# row 2, now shifted to row 4 - the "CNN / Row and Column..." entry) onto the
# two new rows so they pick up the same header-adjacent look (style 1 for
# A:E, style 2 for F:G).
$ws.Range("A4:H4").Copy()
$ws.Range("A2:H2").PasteSpecial(-4122)
$ws.Range("A2:H2").PasteSpecial(-4122)
$ws.Range("A4:H4").PasteSpecial(-4122)

# Row 2: Baseline / intercept only model
$ws.Range("A2").Value = "Baseline"
$ws.Range("B2").Value = "intercept only"
$ws.Range("C2").Value = "Delcode"
$ws.Range("F2").Value = "56.5"

# Row 3: Elastic Net without conn
$ws.Range("A3").Value = "Elastic Net without conn"
$ws.Range("B3").Value = "only age, sex, edyears"
$ws.Range("C3").Value = "Delcode"
$ws.Range("E3").Value = "Jana"
$ws.Range("F3").Value = "71.8"
$ws.Range("G3").Value = "72.6"

# Grow the table so it (and its autofilter) covers the two new rows as well
# as the trailing formatted-but-empty rows that got pushed down.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:H33"))

$ws.Range("J24").Select()
